$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.653.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.16%  "

$ws.Range("D3").Value = "'3.499.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.00%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'589.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.41%  "

$ws.Range("D6").Value = "'186.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.34%  "

$ws.Range("D7").Value = "'0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.46%  "

$ws.Range("D8").Value = "'3.495.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.16%  "

$ws.Range("D10").Value = "'0.175"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("D11").Value = "'0.654"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.39%  "

$ws.Range("D12").Value = "'56.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.54%  "

$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").Value = "'9.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.04%  "

$ws.Range("D15").Value = "'4.045.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.83%  "

$ws.Range("D16").Value = "'18.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.53%  "

$ws.Range("D17").Value = "'3.484.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.72%  "

$ws.Range("D18").Value = "'67.604.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.18%  "

$ws.Range("D19").Value = "'12.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.31%  "

$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("E21").Value = "  +3.78%  "

$ws.Range("D22").Value = "'491.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.37%  "

$ws.Range("D23").Value = "'5.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.81%  "

$ws.Range("D24").Value = "'16.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.00%  "

$ws.Range("D25").Value = "'4.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.59%  "

$ws.Range("D26").Value = "'90.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.47%  "

$ws.Range("D27").Value = "'2.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.08%  "

$ws.Range("D28").Value = "'11.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.06%  "

$ws.Range("D29").Value = "'9.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.94%  "

$ws.Range("D30").Value = "'31.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("D31").Value = "'7.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.54%  "

$ws.Range("D32").Value = "'11.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.15%  "

$ws.Range("D33").Value = "'64.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.40%  "

$ws.Range("D34").Value = "'596.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.27%  "

$ws.Range("E35").Value = "  +4.93%  "

$ws.Range("E36").Value = "  +6.76%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").Value = "'36.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.95%  "

$ws.Range("D39").Value = "'3.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.10%  "

$ws.Range("D40").Value = "'0.388"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.92%  "

$ws.Range("D41").Value = "'0.0₃0771"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.83%  "

$ws.Range("D42").Value = "'3.260.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.37%  "

$ws.Range("D43").Value = "'2.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.78%  "

$ws.Range("D44").Value = "'0.0432"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.88%  "

$ws.Range("E45").Value = "  +3.23%  "

$ws.Range("D46").Value = "'3.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.26%  "

$ws.Range("D47").Value = "'2.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +23.22%  "

$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("D49").Value = "'3.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.29%  "

$ws.Range("D50").Value = "'8.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.01%  "

$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
